# JS-SPA-Self-Evaluation-Protocol.xlsx
# Fill in the "SoftUni Student Info" block (Username / Name / GitHub profile
# link) and turn the GitHub profile link into a real hyperlink, then leave
# the selection on C8 - matching the author's recorded edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Username
$ws.Range("C4").Value = "LittleNinja"

# Name
$ws.Range("C5").Value = "Alexander Stoimenov"

# GitHub profile link - set the visible text first ...
$ws.Range("C7").Value = "https://github.com/astoimenov"

# ... then wire it up as a real (external) hyperlink.
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/astoimenov")

# Leave the selection where the author left it.
$ws.Range("C8").Select()
